$d = $word.ActiveDocument

# Update the date heading
$d.Paragraphs.Item(1).Range.Find.Execute("2023-06-04 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-05 Monday", 2) | Out-Null

# Update each table cell by position (row, col), since some expressions repeat
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "33-8="
$t.Cell(1, 2).Range.Text = "40+44="
$t.Cell(1, 3).Range.Text = "59+4="
$t.Cell(1, 4).Range.Text = "30+46="
$t.Cell(1, 5).Range.Text = "38+12="
$t.Cell(2, 1).Range.Text = "28+8="
$t.Cell(2, 2).Range.Text = "74-2="
$t.Cell(2, 3).Range.Text = "24+22="
$t.Cell(2, 4).Range.Text = "81-36="
$t.Cell(2, 5).Range.Text = "71+15="
$t.Cell(3, 1).Range.Text = "64-2="
$t.Cell(3, 2).Range.Text = "4+28="
$t.Cell(3, 3).Range.Text = "19+55="
$t.Cell(3, 4).Range.Text = "28-20="
$t.Cell(3, 5).Range.Text = "98-10="
$t.Cell(4, 1).Range.Text = "95-35="
$t.Cell(4, 2).Range.Text = "71+16="
$t.Cell(4, 3).Range.Text = "13-13="
$t.Cell(4, 4).Range.Text = "68+20="
$t.Cell(4, 5).Range.Text = "43+35="
$t.Cell(5, 1).Range.Text = "6+91="
$t.Cell(5, 2).Range.Text = "68-47="
$t.Cell(5, 3).Range.Text = "63+18="
$t.Cell(5, 4).Range.Text = "41+46="
$t.Cell(5, 5).Range.Text = "81-42="
$t.Cell(6, 1).Range.Text = "20-18="
$t.Cell(6, 2).Range.Text = "70-13="
$t.Cell(6, 3).Range.Text = "78-46="
$t.Cell(6, 4).Range.Text = "34+24="
$t.Cell(6, 5).Range.Text = "63-35="
$t.Cell(7, 1).Range.Text = "99-10="
$t.Cell(7, 2).Range.Text = "19+15="
$t.Cell(7, 3).Range.Text = "18+41="
$t.Cell(7, 4).Range.Text = "95-37="
$t.Cell(7, 5).Range.Text = "47-27="
$t.Cell(8, 1).Range.Text = "54+0="
$t.Cell(8, 2).Range.Text = "34-30="
$t.Cell(8, 3).Range.Text = "95-59="
$t.Cell(8, 4).Range.Text = "52+10="
$t.Cell(8, 5).Range.Text = "88-62="
$t.Cell(9, 1).Range.Text = "2+80="
$t.Cell(9, 2).Range.Text = "1+79="
$t.Cell(9, 3).Range.Text = "94-41="
$t.Cell(9, 4).Range.Text = "43-2="
$t.Cell(9, 5).Range.Text = "79-37="
$t.Cell(10, 1).Range.Text = "62-55="
$t.Cell(10, 2).Range.Text = "63-34="
$t.Cell(10, 3).Range.Text = "65+1="
$t.Cell(10, 4).Range.Text = "47+29="
$t.Cell(10, 5).Range.Text = "21+68="
$t.Cell(11, 1).Range.Text = "30-30="
$t.Cell(11, 2).Range.Text = "35-30="
$t.Cell(11, 3).Range.Text = "93-51="
$t.Cell(11, 4).Range.Text = "28+57="
$t.Cell(11, 5).Range.Text = "19+40="
$t.Cell(12, 1).Range.Text = "83+10="
$t.Cell(12, 2).Range.Text = "63-14="
$t.Cell(12, 3).Range.Text = "77-71="
$t.Cell(12, 4).Range.Text = "61+13="
$t.Cell(12, 5).Range.Text = "12+82="
$t.Cell(13, 1).Range.Text = "95-39="
$t.Cell(13, 2).Range.Text = "2+41="
$t.Cell(13, 3).Range.Text = "83-51="
$t.Cell(13, 4).Range.Text = "53+32="
$t.Cell(13, 5).Range.Text = "11+9="
$t.Cell(14, 1).Range.Text = "47+17="
$t.Cell(14, 2).Range.Text = "3+43="
$t.Cell(14, 3).Range.Text = "27+67="
$t.Cell(14, 4).Range.Text = "18+12="
$t.Cell(14, 5).Range.Text = "68-52="
$t.Cell(15, 1).Range.Text = "42-41="
$t.Cell(15, 2).Range.Text = "39+30="
$t.Cell(15, 3).Range.Text = "17-9="
$t.Cell(15, 4).Range.Text = "76-8="
$t.Cell(15, 5).Range.Text = "66-62="
$t.Cell(16, 1).Range.Text = "73-49="
$t.Cell(16, 2).Range.Text = "73+23="
$t.Cell(16, 3).Range.Text = "54+38="
$t.Cell(16, 4).Range.Text = "89-37="
$t.Cell(16, 5).Range.Text = "14+9="
$t.Cell(17, 1).Range.Text = "58-41="
$t.Cell(17, 2).Range.Text = "33+28="
$t.Cell(17, 3).Range.Text = "80-60="
$t.Cell(17, 4).Range.Text = "77-58="
$t.Cell(17, 5).Range.Text = "75+1="
$t.Cell(18, 1).Range.Text = "8+10="
$t.Cell(18, 2).Range.Text = "62-16="
$t.Cell(18, 3).Range.Text = "81-72="
$t.Cell(18, 4).Range.Text = "68-15="
$t.Cell(18, 5).Range.Text = "49+37="
$t.Cell(19, 1).Range.Text = "27+10="
$t.Cell(19, 2).Range.Text = "85-35="
$t.Cell(19, 3).Range.Text = "25-3="
$t.Cell(19, 4).Range.Text = "26+34="
$t.Cell(19, 5).Range.Text = "28-7="
$t.Cell(20, 1).Range.Text = "10+69="
$t.Cell(20, 2).Range.Text = "63+32="
$t.Cell(20, 3).Range.Text = "92-69="
$t.Cell(20, 4).Range.Text = "8+14="
$t.Cell(20, 5).Range.Text = "89-31="
